$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (single-dot decimal-looking strings),
# so they stay stored as text just like the rest of column D.
$textCells = @("D5", "D8", "D9", "D11", "D13", "D14", "D19", "D20", "D22", "D27", "D28", "D29", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D42", "D43", "D44", "D45")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "36.531.87"
$ws.Range("E2").Value = "  -1.08%  "

# Row 3
$ws.Range("D3").Value = "2.055.94"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "242.93"
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("E6").Value = "  +1.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "54.53"
$ws.Range("E8").Value = "  -4.63%  "

# Row 9
$ws.Range("D9").Value = "58.26"
$ws.Range("E9").Value = "  -1.80%  "

# Row 10
$ws.Range("E10").Value = "  -3.90%  "

# Row 11
$ws.Range("D11").Value = "0.0749"
$ws.Range("E11").Value = "  -1.86%  "

# Row 12
$ws.Range("E12").Value = "  -3.00%  "

# Row 13
$ws.Range("D13").Value = "0.920"
$ws.Range("E13").Value = "  +5.16%  "

# Row 14
$ws.Range("D14").Value = "14.68"
$ws.Range("E14").Value = "  -4.33%  "

# Row 15
$ws.Range("D15").Value = "2.357.44"
$ws.Range("E15").Value = "  +0.76%  "

# Row 16
$ws.Range("E16").Value = "  -3.70%  "

# Row 17
$ws.Range("D17").Value = "2.056.47"
$ws.Range("E17").Value = "  +2.16%  "

# Row 18
$ws.Range("D18").Value = "36.465.19"
$ws.Range("E18").Value = "  -1.17%  "

# Row 19
$ws.Range("D19").Value = "16.82"
$ws.Range("E19").Value = "  -6.72%  "

# Row 20
$ws.Range("D20").Value = "72.03"
$ws.Range("E20").Value = "  -1.84%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  -2.69%  "

# Row 22
$ws.Range("D22").Value = "238.48"
$ws.Range("E22").Value = "  +1.45%  "

# Row 23
$ws.Range("E23").Value = "  -1.88%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("E25").Value = "  -3.52%  "

# Row 26
$ws.Range("E26").Value = "  -3.09%  "

# Row 27
$ws.Range("D27").Value = "2.13"
$ws.Range("E27").Value = "  +1.45%  "

# Row 28
$ws.Range("D28").Value = "164.22"
$ws.Range("E28").Value = "  -2.84%  "

# Row 29
$ws.Range("D29").Value = "20.06"
$ws.Range("E29").Value = "  +1.31%  "

# Row 30
$ws.Range("E30").Value = "  -0.92%  "

# Row 31
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +11.26%  "

# Row 32
$ws.Range("D32").Value = "5.07"
$ws.Range("E32").Value = "  -5.62%  "

# Row 33
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  -3.32%  "

# Row 34
$ws.Range("E34").Value = "  -2.17%  "

# Row 35
$ws.Range("E35").Value = "  -0.11%  "

# Row 36
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -0.72%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.0826"
$ws.Range("E37").Value = "  -4.68%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("D39").Value = "1.25"
$ws.Range("E39").Value = "  -3.93%  "

# Row 40
$ws.Range("E40").Value = "  -4.72%  "

# Row 41
$ws.Range("E41").Value = "  -2.14%  "

# Row 42
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  -7.54%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.11"
$ws.Range("E43").Value = "  -1.97%  "

# Row 44
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0926"
$ws.Range("E44").Value = "  -4.97%  "

# Row 45
$ws.Range("D45").Value = "93.99"
$ws.Range("E45").Value = "  -2.49%  "

# Row 46
$ws.Range("D46").Value = "1.419.51"
$ws.Range("E46").Value = "  +10.30%  "

# Row 47
$ws.Range("E47").Value = "  +14.01%  "

# Row 48
$ws.Range("E48").Value = "  -5.09%  "

# Row 49
$ws.Range("E49").Value = "  +0.12%  "

# Row 50
$ws.Range("E50").Value = "  -2.34%  "

# Row 51
$ws.Range("D51").Value = "2.247.46"
$ws.Range("E51").Value = "  +0.98%  "

